$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Expertise" (column G) ratings for several rows to reflect
# the new review decisions.
$ws.Range("G6").Value  = "L"
$ws.Range("G8").Value  = ""
$ws.Range("G9").Value  = "L"
$ws.Range("G10").Value = ""
$ws.Range("G13").Value = "L"
$ws.Range("G16").Value = "L"
$ws.Range("G17").Value = ""
$ws.Range("G20").Value = "M"
$ws.Range("G21").Value = "M"
$ws.Range("G23").Value = "L"
$ws.Range("G24").Value = "M"

# Update the active selection to match the saved view state.
$ws.Range("A10:K10").Select()
